$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 (Natalija's record) in place with Aleksandar's data,
# refreshing the timestamp and marking the survey as completed.
$ws.Range("A2").Value = "Aleksandar"
$ws.Range("B2").Value = "Gajic"
$ws.Range("C2").Value = "2023_07_12_10_36_02"
$ws.Range("D2").Value = "gajic7080@gmail.com"
$ws.Range("E2").Value = "Qmi9nT6p5G-tm03YAAsHhCaf_5XHOwXocy9IiDPYSKc"
$ws.Range("F2").Value = $true

# Remove the now-duplicate Aleksandar record that was in row 3.
$ws.Rows("3").Delete()

$wb.Save()
